$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.259.31'
$ws.Range("E2").Value = '  -1.38%  '

$ws.Range("D3").Value = '3.161.60'
$ws.Range("E3").Value = '  +0.79%  '

$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").Value = '''590.96'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.69%  '

$ws.Range("D6").Value = '''138.71'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.49%  '

$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("D8").Value = '3.159.77'
$ws.Range("E8").Value = '  +1.04%  '

$ws.Range("E9").Value = '  -0.72%  '

$ws.Range("E10").Value = '  -1.78%  '

$ws.Range("E11").Value = '  -1.29%  '

$ws.Range("E12").Value = '  -1.70%  '

$ws.Range("E13").Value = '  -2.95%  '

$ws.Range("D14").Value = '''34.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.00%  '

$ws.Range("D15").Value = '3.681.68'
$ws.Range("E15").Value = '  +0.29%  '

$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").Value = '3.157.59'
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = '63.236.09'
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("E19").Value = '  -2.27%  '

$ws.Range("D20").Value = '''477.01'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.76%  '

$ws.Range("D21").Value = '''14.10'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.55%  '

$ws.Range("D22").Value = '''0.703'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("E23").Value = '  +0.95%  '

$ws.Range("D24").Value = '''84.67'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.77%  '

$ws.Range("E25").Value = '  -2.81%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").Value = '''7.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.34%  '

$ws.Range("D29").Value = '''8.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.10%  '

$ws.Range("E30").Value = '  +2.67%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''27.01'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").Value = '''1.00'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.39%  '

$ws.Range("E33").Value = '  -4.50%  '

$ws.Range("E34").Value = '  -4.66%  '

$ws.Range("E35").Value = '  -2.37%  '

$ws.Range("E36").Value = '  -3.38%  '

$ws.Range("D37").Value = '''52.66'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '0.0₃0702'
$ws.Range("E38").Value = '  -7.10%  '

$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("D40").Value = '''422.10'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.95%  '

$ws.Range("E41").Value = '  -8.23%  '

$ws.Range("D42").Value = '''8.30'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.48%  '

$ws.Range("D43").Value = '2.932.66'
$ws.Range("E43").Value = '  +2.21%  '

$ws.Range("E44").Value = '  -5.39%  '

$ws.Range("E45").Value = '  +0.88%  '

$ws.Range("E46").Value = '  -3.83%  '

$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("D48").Value = '''25.55'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.41%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("D50").Value = '''2.26'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -8.35%  '

$ws.Range("D51").Value = '''121.16'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.26%  '
